$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column G ("Recorded By") lists the recording users as a comma-separated
# string. In this edit, whenever "System" appears alongside a real
# recorder account (a "dnasr281@gmail.com" address or
# "backup@backdoor.com"), "System" is swapped with the entry that
# immediately precedes it (e.g. "a, System" -> "System, a", and
# "a, b, System" -> "a, System, b"). Rows that only pair
# "admin@admin.com" with "System", or that don't mix "System" with one
# of those two recorder identities, are left untouched.

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G
    $val = $cell.Value2
    if ([string]::IsNullOrEmpty($val)) { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $rawParts = $val.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    if (-not ($parts -contains "System")) { continue }

    $hasBackdoor = $parts -contains "backup@backdoor.com"
    $hasDnasr = $false
    foreach ($p in $parts) {
        if ($p.ToLower().Contains("dnasr281")) { $hasDnasr = $true }
    }
    if (-not ($hasBackdoor -or $hasDnasr)) { continue }

    $sysIndex = [Array]::IndexOf($parts, "System")
    if ($sysIndex -le 0) { continue }  # nothing precedes it to swap with

    $prevVal = $parts[$sysIndex - 1]
    $parts[$sysIndex - 1] = "System"
    $parts[$sysIndex] = $prevVal

    $newVal = [string]::Join(", ", $parts)
    $ws.Range("G$row").Value = $newVal
}
